# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-level holdings) right before
#    the "总计" (totals) summary sheet.
# 2. Insert a new top row into "总计" summarising the 2022-Q1 holdings,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: style a header cell (bold, centered, top-aligned, thin border)
# ---------------------------------------------------------------------
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# Write a value forcing text ("General"-formatted, unstyled) storage so
# numeric-looking strings (fund codes, scale figures, …) are not coerced
# into numbers - matches the source data convention used on every other
# quarter sheet.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# =======================================================================
# 1) Create the "2022-Q1" worksheet just before "总计"
# =======================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q1.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    Set-HeaderStyle $cell
}

$rows = @(
    @("163415", "兴全商业模式优选混合(LOF)", "159.52", "93.74", "2.95", "4.7058", 9),
    @("001511", "兴全新视野灵活配置定期开放混合", "127.55", "89.43", "2.79", "3.5586", 10),
    @("011738", "华安兴安优选一年持有期混合型证券投资基金A", "25.77", "54.03", "1.49", "0.3840", 4),
    @("011739", "华安兴安优选一年持有期混合型证券投资基金C", "10.03", "54.03", "1.49", "0.1494", 4)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $excelRow = $r + 2
    $rowData = $rows[$r]

    $aCell = $q1.Range("A" + $excelRow)
    $aCell.Value = $r
    Set-HeaderStyle $aCell

    Set-TextValue ($q1.Range("B" + $excelRow)) $rowData[0]
    Set-TextValue ($q1.Range("C" + $excelRow)) $rowData[1]
    Set-TextValue ($q1.Range("D" + $excelRow)) $rowData[2]
    Set-TextValue ($q1.Range("E" + $excelRow)) $rowData[3]
    Set-TextValue ($q1.Range("F" + $excelRow)) $rowData[4]
    Set-TextValue ($q1.Range("G" + $excelRow)) $rowData[5]

    $hCell = $q1.Range("H" + $excelRow)
    $hCell.Value = $rowData[6]
}

# =======================================================================
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet
# =======================================================================
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("B2:D2").ClearFormats()
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 8.800000000000001

# Renumber the 0-based index column (A) for every data row, since the
# pre-existing rows shifted down by one.
for ($r = 2; $r -le 7; $r++) {
    $aCell = $total.Range("A" + $r)
    $aCell.Value = $r - 2
    Set-HeaderStyle $aCell
}

# Restore the original active sheet/selection (sheet creation/edits above
# shift focus onto the sheets we touched last).
$wb.Worksheets.Item(1).Activate()
